# ECLAIRE study-phase source ValueSet workbook:
#  - bump the "Date" metadata cell
#  - split the single "Include from ..." sheet into two:
#      * the existing sheet (rId4/sheetId 2) becomes "Include from ResearchStudyPha"
#        and now points at the HL7 ResearchStudyPhase CodeSystem
#      * a new sheet "Include from Définition des p" is appended right after it,
#        keeping the original ANS/esante CodeSystem content that used to live there
$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item(1)
$metaSheet.Range("B8").Value = "2023-09-01T17:43:23+00:00"

$oldSheet = $wb.Worksheets.Item(2)
$oldName = $oldSheet.Name

# Add the new sheet right after the existing "Include from ..." sheet; it starts
# out blank / default-named, so grab its content + name before anything collides.
$newSheet = $wb.Worksheets.Add($null, $oldSheet)

# Free up the original name on $oldSheet before handing it to $newSheet (sheet
# names must be unique within the workbook).
$oldSheet.Name = "Include from ResearchStudyPha"
$newSheet.Name = $oldName

# Copy the original rows (values + styles + col widths) over to the new sheet
# before we repoint $oldSheet at the new code system.
$oldSheet.Range("A1").Copy($newSheet.Range("A1"))
$oldSheet.Range("A2").Copy($newSheet.Range("A2"))
$oldSheet.Range("A3:B3").Copy($newSheet.Range("A3"))
$oldSheet.Range("A4:B4").Copy($newSheet.Range("A4"))
$newSheet.Columns.Item(1).ColumnWidth = $oldSheet.Columns.Item(1).ColumnWidth
$newSheet.Columns.Item(2).ColumnWidth = $oldSheet.Columns.Item(2).ColumnWidth

# Now update the original sheet's "System URI" value to the new HL7 CodeSystem.
$oldSheet.Range("B4").Value = "http://terminology.hl7.org/CodeSystem/research-study-phase"
